$d = $word.ActiveDocument

# --- Change 1: "...is found active on Indonesia Wikipedia" -> "...is found to be active on Indonesia Wikipedia"
$r1 = $d.Content
$found1 = $r1.Find.Execute("is found active on Indonesia Wikipedia", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "is found to be active on Indonesia Wikipedia", 2)
Write-Host "Change1 found: $found1"

# --- Change 2: " which is a state in India. This particular article lists" -> " a state in India. This particular query lists"
$r2 = $d.Content
$found2 = $r2.Find.Execute("related to Punjab, which is a state in India. This particular article lists out", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "related to Punjab, a state in India. This particular query lists out", 2)
Write-Host "Change2 found: $found2"

# --- Change 3: "By extracting this one can see ... This can be done by sorting the articles by using the category name"
$r3 = $d.Content
$found3 = $r3.Find.Execute("By extracting this one can see the total number of articles that were edited during the competition. This can be done by sorting the articles by using the category name", `
                            $false, $false, $false, $false, $false, $true, 1, $false, `
                            "By extracting this, one can see the total number of articles that were edited during the competition. This is done by sorting the articles using the category name", 2)
Write-Host "Change3 found: $found3"

# --- Change 4: Word Count - 1025 words -> Word Count - 1023 words
$r4 = $d.Content
$found4 = $r4.Find.Execute("Word Count - 1025 words", $false, $false, $false, $false, $false, $true, 1, $false, "Word Count - 1023 words", 2)
Write-Host "Change4 found: $found4"

# --- Change 5: move the _GoBack bookmark from end of doc to within " Kenrick95. Re|trieved October 17"
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$r5 = $d.Content
$found5 = $r5.Find.Execute("Kenrick95. Re", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Change5 found: $found5"
if ($found5) {
    $splitPoint = $r5.End
    $bmRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Host "Done"
